$d = $word.ActiveDocument

# Update the date line in the first paragraph
$d.Content.Find.Execute("2025-02-11 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-02-12 Wednesday", 2)

$table = $d.Tables.Item(1)

# Map of (row, col) -> new text for the 25 populated cells (rows 1,5,9,13,17; cols 1-5)
$updates = @(
    @{Row=1;  Col=1; Text="79÷3=26, 1"},
    @{Row=1;  Col=2; Text="41÷4=10, 1"},
    @{Row=1;  Col=3; Text="74÷2=37, 0"},
    @{Row=1;  Col=4; Text="81÷6=13, 3"},
    @{Row=1;  Col=5; Text="73÷4=18, 1"},

    @{Row=5;  Col=1; Text="70÷4=17, 2"},
    @{Row=5;  Col=2; Text="11÷5=2, 1"},
    @{Row=5;  Col=3; Text="81÷5=16, 1"},
    @{Row=5;  Col=4; Text="86÷3=28, 2"},
    @{Row=5;  Col=5; Text="64÷6=10, 4"},

    @{Row=9;  Col=1; Text="12÷8=1, 4"},
    @{Row=9;  Col=2; Text="13÷9=1, 4"},
    @{Row=9;  Col=3; Text="96÷7=13, 5"},
    @{Row=9;  Col=4; Text="76÷7=10, 6"},
    @{Row=9;  Col=5; Text="86÷6=14, 2"},

    @{Row=13; Col=1; Text="38÷3=12, 2"},
    @{Row=13; Col=2; Text="52÷7=7, 3"},
    @{Row=13; Col=3; Text="27÷3=9, 0"},
    @{Row=13; Col=4; Text="33÷5=6, 3"},
    @{Row=13; Col=5; Text="49÷9=5, 4"},

    @{Row=17; Col=1; Text="79÷3=26, 1"},
    @{Row=17; Col=2; Text="75÷9=8, 3"},
    @{Row=17; Col=3; Text="85÷9=9, 4"},
    @{Row=17; Col=4; Text="62÷6=10, 2"},
    @{Row=17; Col=5; Text="21÷9=2, 3"}
)

foreach ($u in $updates) {
    $cell = $table.Cell($u.Row, $u.Col)
    $cellRange = $cell.Range
    # Trim the trailing cell-mark / paragraph-mark characters from the range
    $cellRange.End = $cellRange.End - 1
    $cellRange.Text = $u.Text
}
